$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (rows 2-5), keeping only the header row
$ws.Range("A2:G5").Clear()

# Update header row: G1 becomes Sub1_Attempt, and add H1:K1 new headers
$ws.Range("G1").Value = "Sub1_Attempt"
$ws.Range("H1").Value = "Sub1_R"
$ws.Range("I1").Value = "Sub1_W"
$ws.Range("J1").Value = "Sub1_NA"
$ws.Range("K1").Value = "Sub1_Marks"
